$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.804914474487305
$ws.Range("B1").Value = 4.424021244049072
$ws.Range("C1").Value = 6.621881008148193
$ws.Range("D1").Value = 7.690051555633545
$ws.Range("E1").Value = 6.094647407531738
